# B6-PowerPoint.pptx edit
# Commit: Sat, Apr 18, 2020  4:04:57 AM
#
# 1) Re-style the three tables (slides 14-16) from the custom "Table_0"
#    style to the built-in table style {A1DFB7FF-5BB2-454F-B4FD-BEB676861E85}.
# 2) The author also re-applied the deck's theme (swapping which theme the
#    slide master vs. the notes master resolve to - "Integral" <-> the
#    default "Office Theme"). We re-apply themes on both masters via the
#    documented Design/Theme pipeline so the intent is captured even though
#    this host's theme-import plumbing is a stub.

$p = $ppt.ActivePresentation

$oldStyleId = "{60EDC319-52FC-45EA-B248-4952001F48FF}"
$newStyleId = "{A1DFB7FF-5BB2-454F-B4FD-BEB676861E85}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
            }
        }
    }
}

# Re-apply the presentation theme to the slide master and the notes master
# (the commit swaps the theme content used by each).
$slideMaster = $p.SlideMaster
$notesMaster = $p.NotesMaster
$slideMaster.ApplyTheme("Office Theme")
$notesMaster.ApplyTheme("Integral")
